$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.808.34'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '3.533.17'
$ws.Range("E3").Value = '  +0.48%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = "'604.45"
$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").Value = "'195.51"
$ws.Range("E6").Value = '  +4.99%  '

$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D9").Value = "'0.203"
$ws.Range("E9").Value = '  -4.91%  '

$ws.Range("D10").Value = "'0.650"
$ws.Range("E10").Value = '  -0.69%  '

$ws.Range("D11").Value = "'53.74"
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("E12").Value = '  -1.15%  '

$ws.Range("D13").Value = "'9.53"
$ws.Range("E13").Value = '  -1.30%  '

$ws.Range("D14").Value = '4.099.31'
$ws.Range("E14").Value = '  +0.34%  '

$ws.Range("D15").Value = "'605.68"
$ws.Range("E15").Value = '  -1.10%  '

$ws.Range("D16").Value = '70.026.43'
$ws.Range("E16").Value = '  +0.31%  '

$ws.Range("D17").Value = "'19.13"
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("D18").Value = "'12.69"
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").Value = '3.539.86'
$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("E20").Value = '  +0.62%  '

$ws.Range("D21").Value = "'0.992"
$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").Value = "'18.18"
$ws.Range("E22").Value = '  +4.01%  '

$ws.Range("D23").Value = "'5.26"
$ws.Range("E23").Value = '  +4.74%  '

$ws.Range("D24").Value = "'102.46"
$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("E25").Value = '  -1.11%  '

$ws.Range("E26").Value = '  +4.52%  '

$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").Value = "'9.62"
$ws.Range("E28").Value = '  -3.12%  '

$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("E30").Value = '  +17.32%  '

$ws.Range("D31").Value = "'7.11"
$ws.Range("E31").Value = '  +1.57%  '

$ws.Range("D32").Value = "'12.58"
$ws.Range("E32").Value = '  +1.43%  '

$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("D34").Value = "'63.20"
$ws.Range("E34").Value = '  -0.69%  '

$ws.Range("D35").Value = '0.0₃0857'
$ws.Range("E35").Value = '  +11.45%  '

$ws.Range("D36").Value = '3.742.06'
$ws.Range("E36").Value = '  +5.38%  '

$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = "'3.03"
$ws.Range("E38").Value = '  -2.69%  '

$ws.Range("D39").Value = "'3.63"
$ws.Range("E39").Value = '  +1.13%  '

$ws.Range("D40").Value = "'0.392"
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("D41").Value = "'36.57"
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").Value = "'488.53"
$ws.Range("E42").Value = '  -7.98%  '

$ws.Range("E43").Value = '  -5.17%  '

$ws.Range("D44").Value = "'0.0457"
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("E45").Value = '  -3.41%  '

$ws.Range("D46").Value = "'0.140"
$ws.Range("E46").Value = '  -1.53%  '

$ws.Range("E47").Value = '  -1.91%  '

$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("D49").Value = "'8.58"
$ws.Range("E49").Value = '  -3.75%  '

$ws.Range("D50").Value = "'0.000253"
$ws.Range("E50").Value = '  +6.16%  '

$ws.Range("D51").Value = "'130.54"
$ws.Range("E51").Value = '  -1.30%  '
